# Insert a new weekly price row for "Vega Modelo de Temuco - Repollo" above the
# existing row 579 (table is kept in date order, so the new record sorts in
# ahead of the previous 2022-01-29 entries), shifting the remaining rows down.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(579).Insert()

$ws.Cells.Item(579, 1).Value2 = 10
$ws.Cells.Item(579, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(579, 3).Value2 = "La Araucanía"
$ws.Cells.Item(579, 4).Value2 = 44585
$ws.Cells.Item(579, 5).Value2 = 9
$ws.Cells.Item(579, 6).Value2 = 100112006
$ws.Cells.Item(579, 7).Value2 = "Repollo"
$ws.Cells.Item(579, 8).Value2 = "Crespo record"
$ws.Cells.Item(579, 9).Value2 = "Primera"
$ws.Cells.Item(579, 10).Value2 = 450
$ws.Cells.Item(579, 11).Value2 = 1000
$ws.Cells.Item(579, 12).Value2 = 1000
$ws.Cells.Item(579, 13).Value2 = 1000
$ws.Cells.Item(579, 14).Value2 = "$/unidad"
$ws.Cells.Item(579, 15).Value2 = "Provincia de Cautín"
$ws.Cells.Item(579, 16).Value2 = 1000
$ws.Cells.Item(579, 17).Value2 = 1
$ws.Cells.Item(579, 18).Value2 = "Hortaliza"
